# "adding averages and more checks"
# Refresh the Training Dashboard's periodic figures (one week later snapshot)
# and flip the now-expired "LOTO (SOPs)" row to NOT VALID; tidy up the
# Exam Dashboard comments column.

$wb = $excel.ActiveWorkbook
$tds = $wb.Worksheets.Item("Training Dashboard")
$exd = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# Training Dashboard: "LAST UPDATE" (col I) moves from 08-Sep-2025 to
# 16-Sep-2025 for every data row, which in turn trims 8 days off the
# "PERIOD TO EXPIRE" (col H) countdown for each row.
# ---------------------------------------------------------------------------

# Column I holds literal text dates ("16-Sep-2025") rather than real Excel
# dates, so pre-format it as Text first -- otherwise entering a date-shaped
# string auto-converts it into a date serial number.
$tds.Range("I3:I34").NumberFormat = "@"

for ($r = 3; $r -le 20; $r++) {
    $tds.Cells.Item($r, 8).Value2 = $tds.Cells.Item($r, 8).Value2 - 8
    $tds.Cells.Item($r, 9).Value2 = "16-Sep-2025"
}

# Row 21 ("LOTO (SOPs)") has now lapsed: 22 days left minus the 8-day jump
# puts it at 14, but its expiry (01-Oct-2025) has been reassessed as already
# past, so the whole row is flipped to the red "NOT VALID" styling (matching
# rows 22/23/25 below it).
$tds.Range("A22:K22").Copy()
$tds.Range("A21:K21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$tds.Cells.Item(21, 8).Value2 = 14
$tds.Cells.Item(21, 9).Value2 = "16-Sep-2025"
$tds.Cells.Item(21, 10).Value2 = "NOT VALID"

for ($r = 22; $r -le 34; $r++) {
    $tds.Cells.Item($r, 8).Value2 = $tds.Cells.Item($r, 8).Value2 - 8
    $tds.Cells.Item($r, 9).Value2 = "16-Sep-2025"
}

# ---------------------------------------------------------------------------
# Header / title styling: bold white text for the title row and the
# dark-blue header row (font no longer carries the old 14pt override).
# ---------------------------------------------------------------------------
foreach ($ws in @($tds, $exd)) {
    $titleRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, 1))
    $titleRange.Font.Size = 11
    $titleRange.Font.Bold = $true
    $titleRange.Font.Color = 16777215

    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, 11))
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = 16777215
}

# ---------------------------------------------------------------------------
# Exam Dashboard: the date-outdated check is rewritten, comments now all
# report "date is valid"; the COMMENTS column is narrowed now that the
# message is shorter.
# ---------------------------------------------------------------------------
$exd.Range("E3").Value2 = "date is valid"
$exd.Range("E4").Value2 = "date is valid"
$exd.Range("E5").Value2 = "date is valid"
$exd.Range("E6").Value2 = "date is valid"

$exd.Columns.Item(5).ColumnWidth = 15
